$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Raw User Data From LDAP"
$ws2 = $wb.Worksheets.Item(2)   # "Workstation Usage Breakdown"

# ---------------------------------------------------------------------------
# Sheet1: "Raw User Data From LDAP"
# ---------------------------------------------------------------------------

# 1. Update "Last updated" date (B1) 2015-08-11 -> 2015-09-21 (serial 40806,
#    date1904 workbook).
$ws1.Range("B1").Value = 40806

# 2. Append the 4 new raw LDAP logins to the bottom of column A.
$ws1.Range("A162").Value = "collinsa"
$ws1.Range("A163").Value = "frechters"
$ws1.Range("A164").Value = "goinac"
$ws1.Range("A165").Value = "ohashi"

# 3. "All Workstation-Only Accounts" table (G/H/I), rows 132-134: fill in the
#    three new people (previously-blank rows). G132 already carries the
#    wrap-text style used throughout the column; G133 loses it; G134 is a
#    brand new row that picks it back up (copied from G132 so no stray style
#    entries get created).
$ws1.Range("G132").Value = "collinsa"
$ws1.Range("H132").Value = "Amanda Collins"
$ws1.Range("I132").Value = "Mouse Light"

$ws1.Range("G133").Value = "frechters"
$ws1.Range("G133").Style = "Normal"
$ws1.Range("H133").Value = "Shahar Frechter"
$ws1.Range("I133").Value = "Visiting Scientist"

$ws1.Range("G134").Value = "ohashi"
$ws1.Range("G132").Copy()
$ws1.Range("G134").PasteSpecial(-4122)
$ws1.Range("H134").Value = "Takako Ohashi"
$ws1.Range("I134").Value = "Mouse Light"

# 4. "All JACS Scientific Accounts" table (C/D/E), rows 138-140: fill in the
#    same three new people. C138 and C140 pick up the column's wrap-text
#    style (copied from an existing styled row), C139 stays plain.
$ws1.Range("C138").Value = "collinsa"
$ws1.Range("C136").Copy()
$ws1.Range("C138").PasteSpecial(-4122)
$ws1.Range("D138").Value = "Amanda Collins"
$ws1.Range("E138").Value = "Mouse Light"

$ws1.Range("C139").Value = "frechters"
$ws1.Range("D139").Value = "Shahar Frechter"
$ws1.Range("E139").Value = "Visiting Scientist"

$ws1.Range("C140").Value = "ohashi"
$ws1.Range("C136").Copy()
$ws1.Range("C140").PasteSpecial(-4122)
$ws1.Range("D140").Value = "Takako Ohashi"
$ws1.Range("E140").Value = "Mouse Light"

$excel.CutCopyMode = $false

# 5. Move the view/selection: no more scrolled-down topLeftCell, cursor on B2.
$ws1.Activate()
$ws1.Range("B2").Select()

# ---------------------------------------------------------------------------
# Sheet2: "Workstation Usage Breakdown"
# ---------------------------------------------------------------------------

# 1. Mouse Light account count 9 -> 10 (the Total at J40 is a SUM formula and
#    recalculates automatically).
$ws2.Range("J26").Value = 10

# 2. Insert a new row for Takako Ohashi right after the other Mouse Light
#    rows (before "chenn", which was row 76) - this shifts every following
#    row down by one.
$ws2.Rows.Item(76).Insert()

$ws2.Range("B76").Value = "ohashi"
$ws2.Range("B76").Font.Color = 0
$ws2.Range("B76").WrapText = $true
$ws2.Range("B76").VerticalAlignment = -4108

$ws2.Range("C76").Value = "Takako Ohashi"
$ws2.Range("C76").Font.Color = 0

$ws2.Range("D76").Value = "Mouse Light"
$ws2.Range("D76").Font.Color = 0

# 3. Clear the stale selection and reset the view to the top of the sheet.
$ws2.Activate()
$ws2.Range("A1").Select()

# Leave "Raw User Data From LDAP" as the active/selected tab, matching the
# target workbook state.
$ws1.Activate()
$ws1.Range("B2").Select()
